$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage for numeric-looking price strings so Excel
# does not silently convert them to floating point numbers (which would drop
# trailing zeros / reformat the text). Number format is reset back to Normal
# right after so no stray cell style is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.277.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4696"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.52%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2871"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06576"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08017"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.50%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.04%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.122"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6847"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "268.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.52%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.264.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007662"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.117.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.270"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.211"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.413"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "

$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.949"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.40%  "

$ws.Range("E29").Value = "  -1.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09863"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.379"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.463"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.074"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.24%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04705"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7004"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.709"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.625"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.295"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.62%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.951"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8426"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4163"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "924.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05676"
$ws.Range("D51").Style = "Normal"

# Row 47 and 48 swap (Aptos <-> EnergySwap)
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.204"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.058"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.20%  "
